$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate stats (row 3-9)
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.49
$summary.Range("B4").Value = 0.5
$summary.Range("B5").Value = 0.13
$summary.Range("B6").Value = 78
$summary.Range("B8").Value = 31
$summary.Range("B9").Value = 41.03

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.49
$status.Range("D4").Value = 78
$status.Range("E4").Value = 0.5
$status.Range("F4").Value = 0.49
$status.Range("G4").Value = 41.03

# ---------------------------------------------------------------------------
# Sheets "All Trades" and "MarketMaking": append new trade #78 as row 79
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(79, 1).Value = 78
    # "2026-02-17" looks like a date to Excel's auto-detection; prefix with an
    # apostrophe (same as a user typing it in) so it is stored as literal text,
    # matching the other Date column cells in this sheet.
    $ws.Cells.Item(79, 2).Value = "'2026-02-17"
    $ws.Cells.Item(79, 3).Value = "08:58:30"
    $ws.Cells.Item(79, 4).Value = "MarketMaking"
    $ws.Cells.Item(79, 5).Value = "UP"
    $ws.Cells.Item(79, 6).Value = 0.14
    $ws.Cells.Item(79, 7).Value = 0.1
    $ws.Cells.Item(79, 8).Value = "CLOSED"
    $ws.Cells.Item(79, 9).Value = -28.5714
    $ws.Cells.Item(79, 10).Value = -0.04
    $ws.Cells.Item(79, 11).Value = 100.49
    $ws.Cells.Item(79, 12).Value = 0
    $ws.Cells.Item(79, 13).Value = 0
    $ws.Cells.Item(79, 14).Value = 0.6
    $ws.Cells.Item(79, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(79, 16).Value = "early_exit"
    $ws.Cells.Item(79, 17).Value = 0.14
}
